# Sprint 3 backlog update:
#   - Task "Assign a Task on a Board from a Project for my Group to an
#     Employee in my Group" (row 10): Week 2 remaining amount goes from 2 to 0
#     (fully tested the controllers -> nothing remaining after week 2).
#   - Task "Assign/Remove Group(s) to/from Stage of Project" (row 11) gets
#     reassigned from Jacob to Daniel.
#   - Leave the selection on C16, matching where the editor left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E10").Value = 0
$ws.Range("B11").Value = "Daniel"

# The SUM formulas in row 14 (and the burndown chart fed by them) recalc
# automatically from these edits.

$ws.Range("C16").Select()
